$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3930.5
$ws.Range("J112").Value = 4925.6665
$ws.Range("L112").Value = 14776.9995
$ws.Range("N112").Value = -16992.9995

$ws.Range("H132").Value = 45459420
$ws.Range("I132").Value = 50005090
$ws.Range("K132").Value = 150015270
$ws.Range("M132").Value = -150012740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2754
$ws.Range("I2").Value = 2692.5
$ws.Range("K2").Value = 2692.5
$ws.Range("M2").Value = -2579.5

$ws.Range("H32").Value = 10340.333
$ws.Range("I32").Value = 9077.799999999999
$ws.Range("K32").Value = 9077.799999999999
$ws.Range("M32").Value = -8790.799999999999

$ws.Range("H45").Value = 1759.6
$ws.Range("I45").Value = 1449.5
$ws.Range("K45").Value = 1449.5
$ws.Range("M45").Value = -1072.5

$ws.Range("H61").Value = 3373.111
$ws.Range("I61").Value = 3373.111
$ws.Range("K61").Value = 3373.111
$ws.Range("M61").Value = -3161.111

$ws.Range("H97").Value = 1086.5333
$ws.Range("I97").Value = 985.6429000000001
$ws.Range("K97").Value = 985.6429000000001
$ws.Range("M97").Value = -489.6429000000001

$ws.Range("H102").Value = 3316.4167
$ws.Range("I102").Value = 3316.4167
$ws.Range("K102").Value = 3316.4167
$ws.Range("M102").Value = -1694.4167

$ws.Range("H110").Value = 2196.8
$ws.Range("I110").Value = 1996.3334
$ws.Range("K110").Value = 1996.3334
$ws.Range("M110").Value = 48.66660000000002

$ws.Range("H116").Value = 2754
$ws.Range("I116").Value = 2692.5
$ws.Range("K116").Value = 2692.5
$ws.Range("M116").Value = -398.5

$ws.Range("H122").Value = 3825.15
$ws.Range("I122").Value = 3926.4736
$ws.Range("K122").Value = 11779.4208
$ws.Range("M122").Value = -9329.4208

$ws.Range("H132").Value = 3348.8484
$ws.Range("I132").Value = 3348.8484
$ws.Range("K132").Value = 10046.5452
$ws.Range("M132").Value = -7516.5452

$ws.Range("H136").Value = 3373.111
$ws.Range("I136").Value = 3373.111
$ws.Range("K136").Value = 10119.333
$ws.Range("M136").Value = -7569.332999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2754
$ws.Range("I3").Value = 2692.5
$ws.Range("K3").Value = 2692.5
$ws.Range("M3").Value = -2578.5

$ws.Range("H99").Value = 2315
$ws.Range("I99").Value = 2315
$ws.Range("K99").Value = 2315
$ws.Range("M99").Value = -817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7387.9287
$ws.Range("I31").Value = 10520.857
$ws.Range("J31").Value = 4255
$ws.Range("K31").Value = 10520.857
$ws.Range("L31").Value = 4255
$ws.Range("M31").Value = -10225.857
$ws.Range("N31").Value = -4845

$ws.Range("H34").Value = 7387.9287
$ws.Range("I34").Value = 10520.857
$ws.Range("J34").Value = 4255
$ws.Range("K34").Value = 10520.857
$ws.Range("L34").Value = 4255
$ws.Range("M34").Value = -10318.857
$ws.Range("N34").Value = -4659

$ws.Range("H122").Value = 2549.75
$ws.Range("I122").Value = 2399.6667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7199.000100000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4749.000100000001
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 9530636
$ws.Range("I132").Value = 10533718
$ws.Range("K132").Value = 31601154
$ws.Range("M132").Value = -31598624

$ws.Range("H134").Value = 2238.2727
$ws.Range("I134").Value = 2228.5264
$ws.Range("K134").Value = 6685.5792
$ws.Range("M134").Value = -4150.5792

$ws.Range("H141").Value = 318180.2
$ws.Range("J141").Value = 318180.2
$ws.Range("L141").Value = 318180.2
$ws.Range("N141").Value = -328540.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1950
$ws.Range("I31").Value = 1950
$ws.Range("K31").Value = 5850
$ws.Range("M31").Value = -5562

$ws.Range("H92").Value = 403.83334
$ws.Range("I92").Value = 396.7
$ws.Range("J92").Value = 439.5
$ws.Range("K92").Value = 1190.1
$ws.Range("L92").Value = 1318.5
$ws.Range("M92").Value = 57.90000000000009
$ws.Range("N92").Value = -3814.5

$ws.Range("H113").Value = 2590.9
$ws.Range("J113").Value = 2618.6
$ws.Range("L113").Value = 7855.799999999999
$ws.Range("N113").Value = -12195.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1621.4828
$ws.Range("I102").Value = 1618.1923
$ws.Range("K102").Value = 1618.1923
$ws.Range("M102").Value = 3.807700000000068

$ws.Range("H122").Value = 2755.25
$ws.Range("I122").Value = 2503
$ws.Range("J122").Value = 3007.5
$ws.Range("K122").Value = 7509
$ws.Range("L122").Value = 9022.5
$ws.Range("M122").Value = -5059
$ws.Range("N122").Value = -13922.5

$ws.Range("H126").Value = 9116
$ws.Range("I126").Value = 8927.857
$ws.Range("K126").Value = 26783.571
$ws.Range("M126").Value = -24313.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1817.3334
$ws.Range("I19").Value = 976.5
$ws.Range("J19").Value = 3499
$ws.Range("K19").Value = 976.5
$ws.Range("L19").Value = 3499
$ws.Range("M19").Value = -806.5
$ws.Range("N19").Value = -3839

$ws.Range("H22").Value = 166667100
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 500000000
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 500000000
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -500000590

$ws.Range("H27").Value = 166667100
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 500000000
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 500000000
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -500000214

$ws.Range("H55").Value = 1271.4
$ws.Range("I55").Value = 674.8
$ws.Range("J55").Value = 1868
$ws.Range("K55").Value = 674.8
$ws.Range("L55").Value = 1868
$ws.Range("M55").Value = -501.8
$ws.Range("N55").Value = -2214

$ws.Range("H132").Value = 4666.6665
$ws.Range("I132").Value = 4666.6665
$ws.Range("K132").Value = 13999.9995
$ws.Range("M132").Value = -11469.9995

$ws.Range("H136").Value = 3716.2727
$ws.Range("I136").Value = 3837.9
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 11513.7
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -8963.700000000001
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 29999.5
$ws.Range("I30").Value = 50000
$ws.Range("J30").Value = 9999
$ws.Range("K30").Value = 50000
$ws.Range("L30").Value = 9999
$ws.Range("M30").Value = -49893
$ws.Range("N30").Value = -10213

$ws.Range("H62").Value = 6500
$ws.Range("I62").Value = 6500
$ws.Range("K62").Value = 6500
$ws.Range("M62").Value = -5876

$ws.Range("H65").Value = 6500
$ws.Range("I65").Value = 6500
$ws.Range("K65").Value = 32500
$ws.Range("M65").Value = -29380

$ws.Range("H107").Value = 437.5
$ws.Range("J107").Value = 403
$ws.Range("L107").Value = 1209
$ws.Range("N107").Value = -5049

$ws.Range("H112").Value = 43554.6
$ws.Range("J112").Value = 43554.6
$ws.Range("L112").Value = 43554.6
$ws.Range("N112").Value = -46508.6

$ws.Range("H136").Value = 6095.6294
$ws.Range("I136").Value = 6784.6523
$ws.Range("J136").Value = 2133.75
$ws.Range("K136").Value = 20353.9569
$ws.Range("L136").Value = 6401.25
$ws.Range("M136").Value = -17803.9569
$ws.Range("N136").Value = -11501.25
